# Fill in the previously-blank English_Reviews_num (G2) and Local_Rank (H2)
# values on the hotel_info sheet of the Houston shard workbook.
#
# The source data keeps these as text-typed cells (they are stored as
# shared strings "1" / "27", not numbers), so we force a text number
# format before assigning the value, then restore the default "Normal"
# style so no stray cell-style index is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "27"
$ws.Range("H2").Style = "Normal"
